# Add new worksheet "total_staff" with staffing data, and make it the active tab.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the existing "annual_budget" sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "total_staff"

# Header row
$ws2.Range("A1").Value = "year"
$ws2.Range("B1").Value = "system"
$ws2.Range("C1").Value = "value"
$ws2.Range("D1").Value = "staff_type"

# Data rows (write the "SUPPORT" string before "SUPERVISION" so the shared
# string table ends up in the same order as the source file)
$ws2.Range("D3").Value = "SUPPORT"
$ws2.Range("D2").Value = "SUPERVISION"

$ws2.Range("A2").Value = 2021
$ws2.Range("B2").Value = "both"
$ws2.Range("C2").Value = 100

$ws2.Range("A3").Value = 2021
$ws2.Range("B3").Value = "both"
$ws2.Range("C3").Value = 50

# Select D2 on the new sheet, and make it the active (visible) sheet/tab.
$ws2.Range("D2").Select()
$ws2.Activate()
